# Commit: "added a few more lines to make it more interesing and louder."
#
# 1) "Oh.." -> "Oh..." and drop the gramStart/gramEnd proofErr wrapper
# 2) "Argg" + "<ellipsis>" -> "*angry grunt*" and drop the spellStart/spellEnd
#    proofErr wrapper
# 3) After the "Why'd I come out here?..." line, add two new lines of
#    dialogue (each preceded by a blank line) and move the _GoBack
#    bookmark down onto the final (still-empty) paragraph.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---- 1) "Oh.. What am I doing here?" paragraph -------------------------
$target = "Oh.. What am I doing here?"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $body = '<w:p><w:pPr><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr>' +
                '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>Oh...</w:t></w:r>' +
                '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> What am I doing here?</w:t></w:r>' +
                '</w:p>'
        $p.Range.InsertXML($pkgHeader + $body + $pkgFooter)
        break
    }
}

# ---- 2) "Argg<ellipsis>" paragraph -> "*angry grunt*" -------------------
$target = "Argg" + [char]8230
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $body = '<w:p><w:pPr><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr>' +
                '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>*angry grunt*</w:t></w:r>' +
                '</w:p>'
        $p.Range.InsertXML($pkgHeader + $body + $pkgFooter)
        break
    }
}

# ---- 3) Grow the ending with two extra lines, relocate the bookmark ----
$rsq = [char]8217  # right single quotation mark, used in the existing text
$target = "Why" + $rsq + "d I come out here? It" + $rsq + "s raining, I don" + $rsq + "t need to water the plants."

$lastIndex = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $lastIndex; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $targetIndex = $i
        break
    }
}

$p1 = $d.Paragraphs($targetIndex)
$p2 = $d.Paragraphs($targetIndex + 1)
$spanStart = $p1.Range.Start
$spanEnd = $p2.Range.End

$pPr = '<w:pPr><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr>'

$para1 = '<w:p>' + $pPr + '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>' + $target + '</w:t></w:r></w:p>'
$paraBlank1 = '<w:p>' + $pPr + '</w:p>'
$paraOh = '<w:p>' + $pPr +
          '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>Oh</w:t></w:r>' +
          '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>,</w:t></w:r>' +
          '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> what' + $rsq + 's that?</w:t></w:r>' +
          '</w:p>'
$paraBlank2 = '<w:p>' + $pPr + '</w:p>'
$paraNoOne = '<w:p>' + $pPr + '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>Oh, no one' + $rsq + 's here.</w:t></w:r></w:p>'
$paraFinal = '<w:p>' + $pPr + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$body = $para1 + $paraBlank1 + $paraOh + $paraBlank2 + $paraNoOne + $paraFinal

$r = $d.Range($spanStart, $spanEnd)
$r.InsertXML($pkgHeader + $body + $pkgFooter)
